# Insert a new weekly record as the first data row of the block (row 292),
# pushing the existing rows 292-412 down to 293-413. Then populate the new
# row 292 with this week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 292..412 down to 293..413 (Excel copies formatting/styles too).
$ws.Rows.Item(292).Insert()

# Fill the new row 292 with the new weekly entry.
$ws.Range("A292").Value = 1
$ws.Range("B292").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C292").Value = "Arica y Parinacota"
$ws.Range("D292").Value = 45134
$ws.Range("E292").Value = 15
$ws.Range("F292").Value = "Fruta"
$ws.Range("G292").Value = 100108
$ws.Range("H292").Value = "Tropicales y subtropicales"
$ws.Range("I292").Value = 100108006
$ws.Range("J292").Value = "Plátano"
$ws.Range("K292").Value = "Sin especificar"
$ws.Range("L292").Value = "Pintón"
$ws.Range("M292").Value = 200
$ws.Range("N292").Value = 15000
$ws.Range("O292").Value = 16000
$ws.Range("P292").Value = 15500
$ws.Range("Q292").Value = "`$/caja 20 kilos"
$ws.Range("R292").Value = "Ecuador"
$ws.Range("S292").Value = 775
$ws.Range("T292").Value = 20
